$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 15243.06
$ws.Range("J17").Value = 15938.828
$ws.Range("L17").Value = 47816.484
$ws.Range("N17").Value = -48152.484
$ws.Range("H76").Value = 10899.5
$ws.Range("I76").Value = 18332.334
$ws.Range("K76").Value = 18332.334
$ws.Range("M76").Value = -18017.334
$ws.Range("H79").Value = 10899.5
$ws.Range("I79").Value = 18332.334
$ws.Range("K79").Value = 18332.334
$ws.Range("M79").Value = -17240.334
$ws.Range("H116").Value = 6251.4287
$ws.Range("J116").Value = 6790.385
$ws.Range("L116").Value = 6790.385
$ws.Range("N116").Value = -13674.385
$ws.Range("H137").Value = 1444756.5
$ws.Range("I137").Value = 1599.4
$ws.Range("J137").Value = 3248702.8
$ws.Range("K137").Value = 4798.200000000001
$ws.Range("L137").Value = 9746108.399999999
$ws.Range("M137").Value = -2248.200000000001
$ws.Range("N137").Value = -9751208.399999999
$ws.Range("H138").Value = 4321.3438
$ws.Range("J138").Value = 7363.857
$ws.Range("L138").Value = 22091.571
$ws.Range("N138").Value = -32371.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1398.8
$ws.Range("I45").Value = 1398.8
$ws.Range("K45").Value = 1398.8
$ws.Range("M45").Value = -1021.8
$ws.Range("H74").Value = 1252199.6
$ws.Range("I74").Value = 1390333
$ws.Range("J74").Value = 8999.5
$ws.Range("K74").Value = 1390333
$ws.Range("L74").Value = 8999.5
$ws.Range("M74").Value = -1389459
$ws.Range("N74").Value = -10747.5
$ws.Range("H77").Value = 1252199.6
$ws.Range("I77").Value = 1390333
$ws.Range("J77").Value = 8999.5
$ws.Range("K77").Value = 6951665
$ws.Range("L77").Value = 44997.5
$ws.Range("M77").Value = -6947297
$ws.Range("N77").Value = -53733.5
$ws.Range("H88").Value = 1914.5385
$ws.Range("I88").Value = 1198.2858
$ws.Range("J88").Value = 2750.1667
$ws.Range("K88").Value = 1198.2858
$ws.Range("L88").Value = 2750.1667
$ws.Range("M88").Value = -792.2858000000001
$ws.Range("N88").Value = -3562.1667
$ws.Range("H91").Value = 1914.5385
$ws.Range("I91").Value = 1198.2858
$ws.Range("J91").Value = 2750.1667
$ws.Range("K91").Value = 1198.2858
$ws.Range("L91").Value = 2750.1667
$ws.Range("M91").Value = 205.7141999999999
$ws.Range("N91").Value = -5558.1667
$ws.Range("H110").Value = 9191.875
$ws.Range("I110").Value = 8922.5
$ws.Range("K110").Value = 8922.5
$ws.Range("M110").Value = -6877.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5265747.5
$ws.Range("I134").Value = 2513.2144
$ws.Range("K134").Value = 7539.6432
$ws.Range("M134").Value = -5004.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5269454
$ws.Range("I16").Value = 7148165
$ws.Range("K16").Value = 7148165
$ws.Range("M16").Value = -7147878
$ws.Range("H31").Value = 31567914
$ws.Range("I31").Value = 37039284
$ws.Range("K31").Value = 37039284
$ws.Range("M31").Value = -37038989
$ws.Range("H34").Value = 31567914
$ws.Range("I34").Value = 37039284
$ws.Range("K34").Value = 37039284
$ws.Range("M34").Value = -37039082
$ws.Range("H107").Value = 4193.7896
$ws.Range("I107").Value = 3825.923
$ws.Range("K107").Value = 3825.923
$ws.Range("M107").Value = -1905.923
$ws.Range("H113").Value = 5269454
$ws.Range("I113").Value = 7148165
$ws.Range("K113").Value = 7148165
$ws.Range("M113").Value = -7145995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 5295.3335
$ws.Range("I121").Value = 747
$ws.Range("K121").Value = 2241
$ws.Range("M121").Value = -931
$ws.Range("H124").Value = 39823.832
$ws.Range("I124").Value = 46343.332
$ws.Range("J124").Value = 33304.332
$ws.Range("K124").Value = 139029.996
$ws.Range("L124").Value = 99912.99600000001
$ws.Range("M124").Value = -134119.996
$ws.Range("N124").Value = -109732.996
$ws.Range("H129").Value = 4999.375
$ws.Range("I129").Value = 3148.6667
$ws.Range("J129").Value = 7378.857
$ws.Range("K129").Value = 9446.000100000001
$ws.Range("L129").Value = 22136.571
$ws.Range("M129").Value = -4446.000100000001
$ws.Range("N129").Value = -32136.571
$ws.Range("H139").Value = 2404.9688
$ws.Range("J139").Value = 3876.8333
$ws.Range("L139").Value = 11630.4999
$ws.Range("N139").Value = -21910.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7925.4165
$ws.Range("I70").Value = 6110.5713
$ws.Range("K70").Value = 6110.5713
$ws.Range("M70").Value = -5840.5713
$ws.Range("H73").Value = 7925.4165
$ws.Range("I73").Value = 6110.5713
$ws.Range("K73").Value = 6110.5713
$ws.Range("M73").Value = -5174.5713
$ws.Range("H97").Value = 1016.73334
$ws.Range("I97").Value = 779.6667
$ws.Range("J97").Value = 1965
$ws.Range("K97").Value = 779.6667
$ws.Range("L97").Value = 1965
$ws.Range("M97").Value = -283.6667
$ws.Range("N97").Value = -2957
$ws.Range("H128").Value = 99999
$ws.Range("J128").Value = 99999
$ws.Range("L128").Value = 99999
$ws.Range("N128").Value = -109959

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 37446.43
$ws.Range("I57").Value = 27021
$ws.Range("K57").Value = 27021
$ws.Range("M57").Value = -26455
$ws.Range("H130").Value = 99998.5
$ws.Range("J130").Value = 99998.5
$ws.Range("L130").Value = 99998.5
$ws.Range("N130").Value = -110038.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 64142.5
$ws.Range("J56").Value = 65000
$ws.Range("L56").Value = 65000
$ws.Range("N56").Value = -66428
$ws.Range("H62").Value = 3040418
$ws.Range("I62").Value = 5350
$ws.Range("J62").Value = 4557952
$ws.Range("K62").Value = 5350
$ws.Range("L62").Value = 4557952
$ws.Range("M62").Value = -4726
$ws.Range("N62").Value = -4559200
$ws.Range("H65").Value = 3040418
$ws.Range("I65").Value = 5350
$ws.Range("J65").Value = 4557952
$ws.Range("K65").Value = 26750
$ws.Range("L65").Value = 22789760
$ws.Range("M65").Value = -23630
$ws.Range("N65").Value = -22796000
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 6675.143
$ws.Range("I126").Value = 6698.316
$ws.Range("K126").Value = 20094.948
$ws.Range("M126").Value = -17624.948
$ws.Range("H133").Value = 59165
$ws.Range("J133").Value = 59165
$ws.Range("L133").Value = 59165
$ws.Range("N133").Value = -69285
$ws.Range("H140").Value = 96332.664
$ws.Range("J140").Value = 96332.664
$ws.Range("L140").Value = 96332.664
$ws.Range("N140").Value = -106692.664
